$d = $word.ActiveDocument

# Helper: build a minimal WordOpenXML "flat OPC" payload wrapping a single
# <w:p> whose children are the run XML fragments supplied, so that
# Range.InsertXML(...) inserts literal <w:r> elements (not just plain text
# that Word would normalize/merge back into the surrounding run).
function New-RunsPayload {
    param([string]$RunsXml)
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p>' + $RunsXml + '</w:p></w:body></w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------
# Change 1: split the "avg_salary" SELECT sentence so that
# "as avg_salary" moves right after AVG(s.salary) instead of after
# t.title, across four separate runs.
# ---------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $full = $p.Range.Text
    if ($full -ne $null -and $full.StartsWith("SELECT e.gender, AVG(s.salary), t.title as avg_salary from")) {
        $start = $p.Range.Start
        $marker = "SELECT e.gender, AVG(s.salary)"
        $cut = $start + $marker.Length
        $tailRange = $d.Range($cut, $start + $full.Length)

        $runsXml = '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
                   '<w:r><w:t>as avg_salary</w:t></w:r>' +
                   '<w:r><w:t>, t.title from employees as e JOIN titles as t on e.emp_id=t.emp_id JOIN salaries as s on e.emp_id=s.emp_id where t.title=&quot;Technique Leader&quot; GROUP by gender</w:t></w:r>'

        $tailRange.InsertXML((New-RunsPayload $runsXml))
        break
    }
}

# ---------------------------------------------------------------------
# Change 2: merge the "dept_emp" run with the following lone-space run
# into a single run "dept_emp ", leaving the surrounding runs (and
# text) untouched.
# ---------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $full = $p.Range.Text
    if ($full -ne $null -and $full.StartsWith("For dept_manager and dept_emp tables:")) {
        $start = $p.Range.Start
        $idx = $full.IndexOf("dept_emp")
        $tailRange = $d.Range($start + $idx, $start + $full.Length)

        $runsXml = '<w:r><w:t xml:space="preserve">dept_emp </w:t></w:r>' +
                   '<w:r><w:t>table</w:t></w:r>' +
                   '<w:r><w:t>s</w:t></w:r>' +
                   '<w:r><w:t xml:space="preserve">: </w:t></w:r>' +
                   '<w:r><w:t xml:space="preserve">There are emp_id and dept_id columns. These columns are likely to repeat the same data, but when </w:t></w:r>' +
                   '<w:r><w:t>combined</w:t></w:r>' +
                   '<w:r><w:t xml:space="preserve"> they form a unique key.</w:t></w:r>'

        $tailRange.InsertXML((New-RunsPayload $runsXml))
        break
    }
}
